# Weekly price update: a new week's reading (2021-10-05) is inserted at the
# top of the "Femacal de La Calera - Ciboulette" series (row 158), pushing
# the existing history (rows 158-183) down by one row (to 159-184).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 158; rows 158..183 shift down to 159..184, carrying
# their formatting (including the date-formatted style on column D) along.
$ws.Rows.Item(158).Insert()

# Populate the newly inserted row 158 with this week's reading. All
# non-date/volume columns repeat the series' static template values.
$ws.Cells.Item(158, 1).Value  = 3
$ws.Cells.Item(158, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(158, 3).Value  = "Coquimbo"
$ws.Cells.Item(158, 4).Value  = 44474
$ws.Cells.Item(158, 5).Value  = 5
$ws.Cells.Item(158, 6).Value  = 100112039
$ws.Cells.Item(158, 7).Value  = "Ciboulette"
$ws.Cells.Item(158, 8).Value  = "Sin especificar"
$ws.Cells.Item(158, 9).Value  = "Primera"
$ws.Cells.Item(158, 10).Value = 160
$ws.Cells.Item(158, 11).Value = 1500
$ws.Cells.Item(158, 12).Value = 1500
$ws.Cells.Item(158, 13).Value = 1500
$ws.Cells.Item(158, 14).Value = "`$/docena de atados"
$ws.Cells.Item(158, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(158, 16).Value = 500
$ws.Cells.Item(158, 17).Value = 3
$ws.Cells.Item(158, 18).Value = "Hortaliza"
